$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: column A label with style copied from A6 (the template data row)
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)
$ws.Range("A7").Value = "a1"

# Row 7: numeric data columns B..CE
$row7vals = @("0", "1.058748410068515e-10", "6.398161279306338e-09", "1.648684159392424e-08", "3.550359217045385e-24", "1.890252129824458e-16", "1.654609060692725e-09", "0.2494240298614347", "0.0001869818372195116", "0.003904829302170598", "0", "0.01662947846207714", "99.69944323157989", "0.01635619339710663", "0.003318085560377239", "0.0007256651528104398", "0.006362514851932103", "0.00180134528534073", "0.001837891659904459", "4.137890119078464e-06", "3.100774848934135e-06", "5.87617926096731e-12", "1.595264144534576e-10", "1.309761518072442e-14", "7.393762874670483e-12", "1.136064747761681e-12", "2.2827027391778e-14", "2.831834781864891e-13", "4.6586782395932e-16", "5.616311077334383e-16", "1.11233340628421e-19", "8.325965345192357e-21", "2.691280473689076e-20", "4.827523493279954e-21", "1.087985574620105e-22", "0", "1.852026186065441e-06", "3.913990061825149e-07", "2.461398264758099e-07", "0", "0", "0", "1.788540095378721e-16", "1.804963092068116e-14", "4.651528904454974e-14", "1.281085076728163e-29", "7.57768665315434e-22", "5.889935027570907e-15", "6.317981587697305e-10", "1.544600554502243e-12", "3.5856195748437e-10", "0", "2.144099294670581e-07", "0.01992665243398281", "3.966321742931323e-05", "0.02889605237384", "0.03229556731351762", "4.041276275273107", "14.20581758498021", "29.53369384332709", "12.51060501543931", "0.05527036064400112", "0.03123957934389144", "0.04815650032836717", "0.03844871584082479", "2.5952886020246", "1.360123230792329", "7.322076507947115", "0.7641682066058152", "18.70049357317203", "0.003336991549543389", "2.984581412190302", "0.134570509131842", "0.008323491422675208", "0.05437399689901418", "2.006542195738793", "0", "2.311725095296766", "0.8458709359497624", "0.3628592253619311", "0", "0")
for ($i = 0; $i -lt $row7vals.Count; $i++) {
    $ws.Cells.Item(7, $i + 2).Value = [double]$row7vals[$i]
}

# Row 8: column A label with style copied from A6 (the template data row)
$ws.Range("A6").Copy()
$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value = "b2"

# Row 8: numeric data columns B..CE
$row8vals = @("0", "1.060301228771484e-10", "6.409693707969709e-09", "1.651596134088414e-08", "3.495735433793975e-24", "1.861762853801977e-16", "1.636339142753496e-09", "0.2494233962652553", "0.0001869773783750931", "0.003904822613890462", "0", "0.01662943896135427", "99.6994439962117", "0.01635619338646314", "0.00331816231728814", "0.0007256596318968826", "0.006362428872269294", "0.001801313301805841", "0.001837856298155758", "4.137814590527763e-06", "3.100613248277515e-06", "5.876129219956474e-12", "1.595171001228785e-10", "1.309757973069216e-14", "7.393731422503394e-12", "1.13608378402341e-12", "2.282772605323726e-14", "2.83185564231592e-13", "4.658740640967222e-16", "5.615973750670056e-16", "1.112357195650552e-19", "8.326054159451756e-21", "2.691086836973063e-20", "4.827239054330904e-21", "1.08800002654838e-22", "0", "1.853459167897335e-06", "3.916823109472563e-07", "2.463499727209743e-07", "0", "0", "0", "1.792966554405978e-16", "1.809924530686002e-14", "4.663537550850302e-14", "1.262222281500869e-29", "7.46649832238714e-22", "5.82793348702725e-15", "6.318035999086779e-10", "1.544582134091942e-12", "3.585667803374677e-10", "0", "2.144180909841319e-07", "0.01992698240228139", "3.966216004171129e-05", "0.02889589380848131", "0.03229558035182906", "4.041276705088187", "14.20581739260818", "29.53369356377717", "12.51060495983317", "0.05527036444642731", "0.03123957926090293", "0.04815650021920672", "0.03844871575306168", "2.595288595901062", "1.360123227655681", "7.32207649141538", "0.7641682048613712", "18.7004935310943", "0.003336991542041314", "2.984581405479622", "0.1345705088292629", "0.008323491403962512", "0.05437399677676604", "2.006542191227598", "0", "2.311725090099648", "0.8458709340481132", "0.3628592245461673", "0", "0")
for ($i = 0; $i -lt $row8vals.Count; $i++) {
    $ws.Cells.Item(8, $i + 2).Value = [double]$row8vals[$i]
}

$excel.CutCopyMode = $false